$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Simoncelli Matteo "
$ws.Range("B22").Value = "Lorenzo Canali | Herta Vernello"
$ws.Range("C22").Value = "Nadir Chtioui | MAI UNA GIOIA"
$ws.Range("D22").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E22").Value = "Michele Merighi | Clitoriders"
$ws.Range("F22").Value = "Matteo  Simoncelli | Herta Vernello"
